# Generate Report for Handoff
# - Update status text "In Translation" -> "Ready for handoff"
# - Bump the handoff timestamps by 40 seconds (12:38:18 -> 12:38:58, 12:38:13 -> 12:38:53)
# - Widen the (now longer) status/date columns to fit the new text

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-28 12:38:58"

$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-28 12:38:53"

$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-28 12:38:58"

$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
